$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for columns B..G (column A header already exists)
$ws.Range("B1").Value = "aspect ratio"
$ws.Range("C1").Value = "extent"
$ws.Range("D1").Value = "Blue"
$ws.Range("E1").Value = "Green"
$ws.Range("F1").Value = "Red"
$ws.Range("G1").Value = "Hue"

# Copy style (bold font, thin border, centered/top alignment) from A1 to the new header cells
$ws.Range("A1").Copy()
$ws.Range("B1:G1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Replace the numeric data (rows 2-21, columns A-G)
$ws.Range("A2").Value = 65
$ws.Range("B2").Value = 0.1719789131198317
$ws.Range("C2").Value = 0.8982768326077035
$ws.Range("D2").Value = 29.20517301642782
$ws.Range("E2").Value = 87.63998601887452
$ws.Range("F2").Value = 111.1006641034603
$ws.Range("G2").Value = 15.91261796574624

$ws.Range("A3").Value = 69
$ws.Range("B3").Value = 0.1827491780715877
$ws.Range("C3").Value = 0.8435412556136059
$ws.Range("D3").Value = 34.04604726976365
$ws.Range("E3").Value = 89.2322738386308
$ws.Range("F3").Value = 113.4099429502852
$ws.Range("G3").Value = 14.28198859005705

$ws.Range("A4").Value = 55
$ws.Range("B4").Value = 0.2001229539807216
$ws.Range("C4").Value = 0.8239227184043145
$ws.Range("D4").Value = 36.64510412051396
$ws.Range("E4").Value = 89.1856446610545
$ws.Range("F4").Value = 111.1789986708019
$ws.Range("G4").Value = 18.83473637571998

$ws.Range("A5").Value = 57
$ws.Range("B5").Value = 0.1724024047621364
$ws.Range("C5").Value = 0.8603092840105619
$ws.Range("D5").Value = 29.96866230121609
$ws.Range("E5").Value = 83.80308699719365
$ws.Range("F5").Value = 106.6721234798877
$ws.Range("G5").Value = 18.93592142188962

$ws.Range("A6").Value = 56
$ws.Range("B6").Value = 0.1709531994201213
$ws.Range("C6").Value = 0.8511780904461099
$ws.Range("D6").Value = 28.61282771535581
$ws.Range("E6").Value = 88.08192883895131
$ws.Range("F6").Value = 106.3909176029963
$ws.Range("G6").Value = 16.80571161048689

$ws.Range("A7").Value = 52
$ws.Range("B7").Value = 0.1705332642443037
$ws.Range("C7").Value = 0.8710008561078881
$ws.Range("D7").Value = 34.12296650717703
$ws.Range("E7").Value = 90.38612440191388
$ws.Range("F7").Value = 111.6253588516746
$ws.Range("G7").Value = 15.99760765550239

$ws.Range("A8").Value = 54
$ws.Range("B8").Value = 0.1765076877759789
$ws.Range("C8").Value = 0.8434846771593958
$ws.Range("D8").Value = 28.2516156828953
$ws.Range("E8").Value = 89.32701421800948
$ws.Range("F8").Value = 112.2856527358897
$ws.Range("G8").Value = 17.57690650581646

$ws.Range("A9").Value = 61
$ws.Range("B9").Value = 0.1645684068006661
$ws.Range("C9").Value = 0.8994874721960268
$ws.Range("D9").Value = 28.53610988037218
$ws.Range("E9").Value = 82.89455028799291
$ws.Range("F9").Value = 108.7771377935312
$ws.Range("G9").Value = 16.07665042091272

$ws.Range("A10").Value = 45
$ws.Range("B10").Value = 0.1660514160824168
$ws.Range("C10").Value = 0.8758838208671624
$ws.Range("D10").Value = 23.96145610278372
$ws.Range("E10").Value = 88.77141327623126
$ws.Range("F10").Value = 100.4544967880086
$ws.Range("G10").Value = 21.63222698072805

$ws.Range("A11").Value = 39
$ws.Range("B11").Value = 0.1730335239697859
$ws.Range("C11").Value = 0.8412769322230743
$ws.Range("D11").Value = 32.61223431829963
$ws.Range("E11").Value = 88.39450492483152
$ws.Range("F11").Value = 110.4271643338517
$ws.Range("G11").Value = 17.23224468636599

$ws.Range("A12").Value = 40
$ws.Range("B12").Value = 0.1782651100834997
$ws.Range("C12").Value = 0.8283141040512229
$ws.Range("D12").Value = 31.91084462982273
$ws.Range("E12").Value = 86.67987486965589
$ws.Range("F12").Value = 113.633472367049
$ws.Range("G12").Value = 21.82012513034411

$ws.Range("A13").Value = 44
$ws.Range("B13").Value = 0.169620890396565
$ws.Range("C13").Value = 0.8989907191017944
$ws.Range("D13").Value = 29.5310119695321
$ws.Range("E13").Value = 90.93362350380849
$ws.Range("F13").Value = 109.9646354733406
$ws.Range("G13").Value = 21.84058759521219

$ws.Range("A14").Value = 41
$ws.Range("B14").Value = 0.173251042036921
$ws.Range("C14").Value = 0.8441967220738036
$ws.Range("D14").Value = 31.04685890834191
$ws.Range("E14").Value = 86.90061791967044
$ws.Range("F14").Value = 112.711122554068
$ws.Range("G14").Value = 20.78836251287333

$ws.Range("A15").Value = 59
$ws.Range("B15").Value = 0.1663780661615406
$ws.Range("C15").Value = 0.8809718603253573
$ws.Range("D15").Value = 33.56353839245622
$ws.Range("E15").Value = 86.79928154467893
$ws.Range("F15").Value = 115.3632689717108
$ws.Range("G15").Value = 20.40817242927706

$ws.Range("A16").Value = 45
$ws.Range("B16").Value = 0.1658502492317681
$ws.Range("C16").Value = 0.8808121402609804
$ws.Range("D16").Value = 27.73424796747967
$ws.Range("E16").Value = 88.2428861788618
$ws.Range("F16").Value = 107.8272357723577
$ws.Range("G16").Value = 13.59146341463415

$ws.Range("A17").Value = 74
$ws.Range("B17").Value = 0.1685206107547761
$ws.Range("C17").Value = 0.8781752109968816
$ws.Range("D17").Value = 28.72805507745267
$ws.Range("E17").Value = 90.86617900172116
$ws.Range("F17").Value = 122.1652323580034
$ws.Range("G17").Value = 11.87177280550775

$ws.Range("A18").Value = 74
$ws.Range("B18").Value = 0.1834506928529379
$ws.Range("C18").Value = 0.8466548770919182
$ws.Range("D18").Value = 31.6819801980198
$ws.Range("E18").Value = 86.70653465346534
$ws.Range("F18").Value = 120.8641584158416
$ws.Range("G18").Value = 11.74178217821782

$ws.Range("A19").Value = 78
$ws.Range("B19").Value = 0.1767236012286736
$ws.Range("C19").Value = 0.878897824486619
$ws.Range("D19").Value = 30.18401639344262
$ws.Range("E19").Value = 89.72868852459017
$ws.Range("F19").Value = 117.5139344262295
$ws.Range("G19").Value = 12.0405737704918

$ws.Range("A20").Value = 70
$ws.Range("B20").Value = 0.1662980987940527
$ws.Range("C20").Value = 0.941033328756802
$ws.Range("D20").Value = 32.89067974772249
$ws.Range("E20").Value = 88.54204625087597
$ws.Range("F20").Value = 123.177295024527
$ws.Range("G20").Value = 13.89453398738612

$ws.Range("A21").Value = 71
$ws.Range("B21").Value = 0.164824341152918
$ws.Range("C21").Value = 0.8927277137444772
$ws.Range("D21").Value = 27.06660666066607
$ws.Range("E21").Value = 86.95679567956796
$ws.Range("F21").Value = 118.3532853285328
$ws.Range("G21").Value = 10.76237623762376
